$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the old second header row ("Hiver / Été / Année") so that the
#    13 data rows shift up by one (old rows 3-15 become new rows 2-14).
$ws.Rows(2).Delete()

# 2. Build the new combined header row (row 1).
#    Columns A-E are the new identification/date columns (default style).
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# 3. Create a lightweight named style that only carries the data font
#    (Arial 9, same as the existing font used across the sheet) with no
#    explicit number format, then apply it to the measurement headers
#    F1:K1. Deleting the temporary named style afterwards leaves its
#    cell-format entry behind (still referenced by F1:K1) without
#    polluting the workbook with an extra named cell style.
$headerStyle = $wb.Styles.Add("NE2008HeaderStyle")
$headerStyle.Font.Name = "Arial"
$headerStyle.Font.Size = 9

$ws.Range("F1:K1").Style = "NE2008HeaderStyle"
$wb.Styles.Item("NE2008HeaderStyle").Delete()

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# 4. Match the refreshed selection highlighted in the saved workbook.
$ws.Range("A2:K2").Select()
